$wb = $excel.ActiveWorkbook

# OFF sheet - Row 3 ("R") updates
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 273
$wsOff.Range("C3").Value = 185
$wsOff.Range("D3").Value = 69
$wsOff.Range("E3").Value = 41
$wsOff.Range("F3").Value = 6

# DEF sheet - Row 3 ("R") updates
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 290
$wsDef.Range("C3").Value = 216
$wsDef.Range("D3").Value = 66
$wsDef.Range("E3").Value = 31
$wsDef.Range("F3").Value = 3
